$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the hyperlinks that live on row 4 (C4:G4) before clearing the
# cells, since clearing contents alone does not drop the hyperlink
# objects / their relationships. Re-query the (live) Hyperlinks
# collection for every single deletion so stale indices/object
# references from a previous deletion can't shift what gets removed.
function Remove-HyperlinkAt($addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            $hl.Delete()
            return
        }
    }
}

Remove-HyperlinkAt '$C$4'
Remove-HyperlinkAt '$D$4'
Remove-HyperlinkAt '$E$4'
Remove-HyperlinkAt '$F$4'
Remove-HyperlinkAt '$G$4'

# Clear the data row (row 4): A4 and C4:G4 keep their style but lose
# their values, while B4 is wiped completely (both value and format)
# so the cell disappears from the sheet entirely.
$ws.Range("A4").ClearContents()
$ws.Range("B4").Clear()
$ws.Range("C4:G4").ClearContents()

# Update the selection to match the new state (whole row 4 selected).
$null = $ws.Rows.Item(4).Select()
